$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated B2:B118 accuracy values
$values = @(0.9375, 0.875, 0.84375, 0.78125, 0.703125, 0.6875, 0.65625, 0.640625, 0.578125, 0.609375, 0.625, 0.59375, 0.609375, 0.59375, 0.59375, 0.625, 0.46875, 0.5625, 0.484375, 0.5, 0.53125, 0.484375, 0.484375, 0.484375, 0.484375, 0.46875, 0.5, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.46875, 0.453125, 0.421875, 0.421875, 0.421875, 0.421875, 0.421875, 0.4375, 0.4375, 0.421875, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.4375, 0.453125, 0.453125, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.421875, 0.578125, 0.453125, 0.5, 0.515625, 0.5, 0.5, 0.53125, 0.5, 0.5, 0.546875, 0.53125, 0.484375, 0.5, 0.546875, 0.4098360655737705)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Rows 102-118 column A text changed (Python repr object memory address)
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = "<__main__.DisplayOutputs object at 0x7f00e84e7520>"
}

# Select the whole sheet (Ctrl+A), active cell remains within the selection
$ws.Cells.Select()
